# October 2018 dashboard update: append a new "October 2018" data row (row 19)
# to Sheet1, mirroring the existing monthly rows (2-18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A19: month label "October\n2018". Typed directly, Excel's smart-entry
# would parse this as a date (it does even with a text number format), so we
# build it as a tiny text formula in a scratch cell, paste its *value* into
# A19 (landing as a real shared-string literal), then copy the format of an
# existing month cell (A9) onto it so it keeps the mmm/yyyy + wrap-text look
# used by every other month header in column A. ---
$monthLabel = '="October' + [char]10 + '2018"'
$ws.Range("Z1").Formula = $monthLabel
$ws.Range("Z1").Copy()
$ws.Range("A19").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("Z1").Clear()
$ws.Range("A9").Copy()
$ws.Range("A19").PasteSpecial(-4122)   # xlPasteFormats

# --- Numeric inputs (columns B:J) ---
$ws.Range("B19").Value = 303000
$ws.Range("C19").Value = 383000
$ws.Range("D19").Value = 91500
$ws.Range("E19").Value = 11500
$ws.Range("F19").Value = 33000
$ws.Range("G19").Value = 519000
$ws.Range("H19").Value = 110000
$ws.Range("I19").Value = 118000
$ws.Range("J19").Value = 3000

# K19 stays blank (formatted, like the rest of the column)
$ws.Range("K19").Value = $null

# L19: percentage input
$ws.Range("L19").Value = 0.59

# M19 stays blank
$ws.Range("M19").Value = $null

# --- Formulas (columns N:Q), same shape as row 18 ---
$ws.Range("N19").Formula = "=IF(H19>I19,C19-(H19-I19),C19)*-1-R19"
$ws.Range("O19").Formula = "=IF(I19>H19,I19-H19,0)"
$ws.Range("P19").Formula = "=IF(H19>I19,H19-I19,0)*-1"
$ws.Range("Q19").Formula = "=IF(H19>I19,I19,H19)"

# R19 stays blank
$ws.Range("R19").Value = $null

# --- S:U, reuse the existing column styles/number formats exactly ---
$ws.Range("S18:U18").Copy()
$ws.Range("S19").PasteSpecial(-4122)
$ws.Range("S19").Value = 1626.8
$ws.Range("T19").Value = $null
$ws.Range("U19").Value = 80000

# --- Give the new row's B:R cells the distinct font used for this dataset
# update, with the matching number formats (thousands separator / percent /
# general). M19 (plain General + new font) first, so it lands on the same
# "base" new style the rest of the row's new styles build on. ---
$ws.Range("M19").Font.Color = 0

$ws.Range("B19:J19").NumberFormat = "#,##0"
$ws.Range("B19:J19").Font.Color = 0

$ws.Range("L19").NumberFormat = "0%"
$ws.Range("L19").Font.Color = 0

$ws.Range("K19").NumberFormat = "#,##0"
$ws.Range("K19").Font.Color = 0

$ws.Range("N19:Q19").NumberFormat = "#,##0"
$ws.Range("N19:Q19").Font.Color = 0

$ws.Range("R19").NumberFormat = "#,##0"
$ws.Range("R19").Font.Color = 0

# Row height matches the other "ht=30" rows (e.g. row 17)
$ws.Rows.Item(19).RowHeight = 30

# --- View state: scroll so row ~10 is near the top, select T21 (matches the
# saved sheetView in the target workbook) ---
$ws.Range("T21").Select()
